# Updated cryptos list on Sun Jul 16 17:53:08 UTC 2023 with GitHub Actions
#
# Refresh the scraped coinranking.com price/volume snapshot. The "Price"
# column (D) holds values that look numeric (e.g. "0.3210", "30.361.71")
# but must stay plain text exactly as scraped (leading/trailing zeros,
# thousand-dot-separated big prices, etc.), matching the original
# inline-string cells. Force text via NumberFormat "@" before writing,
# then restore the default "Normal" style so no stray formatting is left
# on the cell (it only had the default style to begin with).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = "D2"; Value = "30.361.71" },
    @{ Ref = "E2"; Value = "  +0.16%  " },
    @{ Ref = "D3"; Value = "1.935.69" },
    @{ Ref = "E3"; Value = "  +0.16%  " },
    @{ Ref = "E4"; Value = "  -0.25%  " },
    @{ Ref = "D5"; Value = "0.7733" },
    @{ Ref = "E5"; Value = "  +8.66%  " },
    @{ Ref = "D6"; Value = "246.49" },
    @{ Ref = "E6"; Value = "  -1.69%  " },
    @{ Ref = "D7"; Value = "0.9989" },
    @{ Ref = "E7"; Value = "  -0.17%  " },
    @{ Ref = "D8"; Value = "27.99" },
    @{ Ref = "E8"; Value = "  +1.74%  " },
    @{ Ref = "D9"; Value = "0.3210" },
    @{ Ref = "E9"; Value = "  -2.57%  " },
    @{ Ref = "D10"; Value = "0.07089" },
    @{ Ref = "E10"; Value = "  -2.07%  " },
    @{ Ref = "D11"; Value = "0.7826" },
    @{ Ref = "E11"; Value = "  -2.62%  " },
    @{ Ref = "E12"; Value = "  -0.82%  " },
    @{ Ref = "D13"; Value = "1.934.11" },
    @{ Ref = "E13"; Value = "  +0.10%  " },
    @{ Ref = "D14"; Value = "5.378" },
    @{ Ref = "E14"; Value = "  -1.51%  " },
    @{ Ref = "D15"; Value = "94.97" },
    @{ Ref = "E15"; Value = "  +0.61%  " },
    @{ Ref = "D16"; Value = "14.55" },
    @{ Ref = "E16"; Value = "  -3.22%  " },
    @{ Ref = "D17"; Value = "30.360.57" },
    @{ Ref = "D18"; Value = "256.08" },
    @{ Ref = "E18"; Value = "  +1.46%  " },
    @{ Ref = "D19"; Value = "0.000008010" },
    @{ Ref = "E19"; Value = "  -1.85%  " },
    @{ Ref = "D20"; Value = "5.824" },
    @{ Ref = "E20"; Value = "  +0.72%  " },
    @{ Ref = "D21"; Value = "2.187.91" },
    @{ Ref = "E21"; Value = "  +0.07%  " },
    @{ Ref = "E22"; Value = "  -0.12%  " },
    @{ Ref = "D23"; Value = "0.9978" },
    @{ Ref = "E23"; Value = "  -0.43%  " },
    @{ Ref = "D24"; Value = "6.750" },
    @{ Ref = "E24"; Value = "  -3.08%  " },
    @{ Ref = "D25"; Value = "9.604" },
    @{ Ref = "E25"; Value = "  -1.24%  " },
    @{ Ref = "D26"; Value = "164.03" },
    @{ Ref = "E26"; Value = "  -0.76%  " },
    @{ Ref = "D27"; Value = "0.1349" },
    @{ Ref = "E27"; Value = "  +4.90%  " },
    @{ Ref = "D28"; Value = "19.11" },
    @{ Ref = "E28"; Value = "  -0.83%  " },
    @{ Ref = "D29"; Value = "2.290" },
    @{ Ref = "E29"; Value = "  -2.38%  " },
    @{ Ref = "E30"; Value = "  +1.11%  " },
    @{ Ref = "D31"; Value = "1.519" },
    @{ Ref = "E31"; Value = "  -1.51%  " },
    @{ Ref = "D32"; Value = "4.429" },
    @{ Ref = "E32"; Value = "  +0.33%  " },
    @{ Ref = "D33"; Value = "4.146" },
    @{ Ref = "E33"; Value = "  -0.51%  " },
    @{ Ref = "D34"; Value = "0.05191" },
    @{ Ref = "E34"; Value = "  +0.08%  " },
    @{ Ref = "D35"; Value = "1.284" },
    @{ Ref = "E35"; Value = "  +2.10%  " },
    @{ Ref = "D36"; Value = "0.7527" },
    @{ Ref = "E36"; Value = "  +1.09%  " },
    @{ Ref = "D37"; Value = "2.771" },
    @{ Ref = "E37"; Value = "  -0.76%  " },
    @{ Ref = "D38"; Value = "0.01975" },
    @{ Ref = "E38"; Value = "  +0.63%  " },
    @{ Ref = "D39"; Value = "2.814" },
    @{ Ref = "E39"; Value = "  +0.21%  " },
    @{ Ref = "D40"; Value = "78.99" },
    @{ Ref = "E40"; Value = "  +0.49%  " },
    @{ Ref = "D41"; Value = "6.448" },
    @{ Ref = "E41"; Value = "  +0.53%  " },
    @{ Ref = "D42"; Value = "0.4520" },
    @{ Ref = "E42"; Value = "  +0.14%  " },
    @{ Ref = "D43"; Value = "1.980" },
    @{ Ref = "E43"; Value = "  -1.61%  " },
    @{ Ref = "E44"; Value = "  -0.06%  " },
    @{ Ref = "D45"; Value = "0.8344" },
    @{ Ref = "E45"; Value = "  -1.10%  " },
    @{ Ref = "D46"; Value = "101.20" },
    @{ Ref = "E46"; Value = "  -0.12%  " },
    @{ Ref = "D47"; Value = "9.787" },
    @{ Ref = "E47"; Value = "  +0.50%  " },
    @{ Ref = "D48"; Value = "7.509" },
    @{ Ref = "E48"; Value = "  +1.13%  " },
    @{ Ref = "D49"; Value = "37.44" },
    @{ Ref = "E49"; Value = "  +2.13%  " },
    @{ Ref = "D50"; Value = "982.99" },
    @{ Ref = "E50"; Value = "  +11.24%  " },
    @{ Ref = "B51"; Value = "Decentraland" },
    @{ Ref = "C51"; Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana" },
    @{ Ref = "D51"; Value = "0.4168" },
    @{ Ref = "E51"; Value = "  +0.20%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    if ($u.Ref -like "D*") {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
